$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "SIGNAL VALID NAME" column values for rows 2-5
# (order chosen to reproduce the original shared-string insertion order)
$ws.Range("J2").Value = "AXICRYPT_AXI_MI0_ARVALID"
$ws.Range("J5").Value = "AXICRYPT_AXI_MI0_WVALID"
$ws.Range("J4").Value = "AXICRYPT_AXI_MI0_RVALID"
$ws.Range("J3").Value = "AXICRYPT_AXI_MI0_AWVALID"

# Fix GROUP column value for row 4 (was "B", now "A")
$ws.Range("H4").Value = "A"

# Resize column J to fit the new, wider content (matches Excel's own best-fit result of 26 characters)
$ws.Columns.Item(10).ColumnWidth = 25.2

# Update the active selection to reflect the last edited cell
$ws.Range("J10").Select() | Out-Null
